$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# Add the new worksheet right after "login"
$ws = $wb.Worksheets.Add($null, $loginSheet)
$ws.Name = "user_details"

# Populate header row (order mirrors the original author's entry order)
$ws.Range("A1").Value = "Firstname"
$ws.Range("C1").Value = "Postalcode"
$ws.Range("B1").Value = "Lastname"

# Populate data row
$ws.Range("A2").Value = "Selina"
$ws.Range("B2").Value = "Mabunda"
$ws.Range("C2").Value = 1680

# Reuse the bordered style already used on the "login" sheet's data rows
$loginSheet.Range("A2").Copy()
$null = $ws.Range("A1:C4").PasteSpecial(-4122)

# Column C width
$ws.Columns("C").ColumnWidth = 9.1666666666666667

# Selection matches the committed sheet (C2 selected)
$null = $ws.Range("C2").Select()

Write-Host "done"
